$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$ws.Cells.Item(2, 7).Value = 'shaimaa.ahmed@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
$ws.Cells.Item(4, 7).Value = 'hend_mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, heba@med.asu.edu.eg'
$ws.Cells.Item(5, 7).Value = 'Fatmaelhady@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Nadagouda@med.asu.edu.eg'
$ws.Cells.Item(8, 7).Value = 'aml.awwad@med.asu.edu.eg, aya.saeed@med.asu.edu.eg'
$ws.Cells.Item(10, 7).Value = 'Madeha.Saeed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg'
$ws.Cells.Item(11, 7).Value = 'Sarah.Mahdy@med.asu.edu.eg, Alaa.A.Mostafa@med.asu.edu.eg'
$ws.Cells.Item(12, 7).Value = 'Sarah.Mahdy@med.asu.edu.eg, Noran.Mahmoud@med.asu.edu.eg'
$ws.Cells.Item(18, 7).Value = 'yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg'
$ws.Cells.Item(19, 7).Value = 'yassmen.ahmed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Cells.Item(20, 7).Value = 'yassmen.ahmed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, marinasorial@med.asu.edu.eg'
$ws.Cells.Item(21, 7).Value = 'yasmin.m.senosy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg'
$ws.Cells.Item(22, 7).Value = 'Monica.Eshak@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg'
$ws.Cells.Item(23, 7).Value = 'yassmen.ahmed@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg'
$ws.Cells.Item(24, 7).Value = 'shaimaa.ahmed@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
$ws.Cells.Item(25, 7).Value = 'manar.montaser@med.asu.edu.eg, gehanadel@med.asu.edu.eg, backup@backdoor.com, alshimaa.atef@med.asu.edu.egm'
$ws.Cells.Item(26, 7).Value = 'hend_mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, heba@med.asu.edu.eg'
$ws.Cells.Item(27, 7).Value = 'Fatmaelhady@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Nadagouda@med.asu.edu.eg'
$ws.Cells.Item(30, 7).Value = 'aml.awwad@med.asu.edu.eg, aya.saeed@med.asu.edu.eg'
$ws.Cells.Item(32, 7).Value = 'Madeha.Saeed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg'
$ws.Cells.Item(33, 7).Value = 'Sarah.Mahdy@med.asu.edu.eg, Alaa.A.Mostafa@med.asu.edu.eg'
$ws.Cells.Item(34, 7).Value = 'Sarah.Mahdy@med.asu.edu.eg, Noran.Mahmoud@med.asu.edu.eg'
$ws.Cells.Item(40, 7).Value = 'yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg'
$ws.Cells.Item(41, 7).Value = 'yassmen.ahmed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Cells.Item(42, 7).Value = 'yassmen.ahmed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, marinasorial@med.asu.edu.eg'
$ws.Cells.Item(43, 7).Value = 'yasmin.m.senosy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg'
$ws.Cells.Item(44, 7).Value = 'Monica.Eshak@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg'
$ws.Cells.Item(45, 7).Value = 'yassmen.ahmed@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg'
$ws.Cells.Item(46, 7).Value = 'hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg'
$ws.Cells.Item(48, 7).Value = 'nahla.nagiub@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, heba@med.asu.edu.eg'
$ws.Cells.Item(49, 7).Value = 'menna-alah.mohamed@asu.edu.eg, Nadagouda@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg'
$ws.Cells.Item(52, 7).Value = 'Shimaa.ashraf@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg'
$ws.Cells.Item(54, 7).Value = 'Madeha.Saeed@med.asu.edu.eg, merna.said@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, maimustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg'
$ws.Cells.Item(60, 7).Value = 'user@user.com, nancy.abdelshafy@med.asu.edu.eg, Amr-Saeed@med.asu.edu.eg'
$ws.Cells.Item(62, 7).Value = 'yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg'
$ws.Cells.Item(63, 7).Value = 'Monica.Eshak@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg'
$ws.Cells.Item(64, 7).Value = 'youstina.magdy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Cells.Item(65, 7).Value = 'eman.samir@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Cells.Item(66, 7).Value = 'maryam.ashraf@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Eman_mohamed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, marinasorial@med.asu.edu.eg'
$ws.Cells.Item(68, 7).Value = 'hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg'
$ws.Cells.Item(70, 7).Value = 'nahla.nagiub@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, heba@med.asu.edu.eg'
$ws.Cells.Item(71, 7).Value = 'menna-alah.mohamed@asu.edu.eg, Nadagouda@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg'
$ws.Cells.Item(74, 7).Value = 'Shimaa.ashraf@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg'
$ws.Cells.Item(76, 7).Value = 'Madeha.Saeed@med.asu.edu.eg, merna.said@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, maimustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg'
$ws.Cells.Item(82, 7).Value = 'user@user.com, nancy.abdelshafy@med.asu.edu.eg, Amr-Saeed@med.asu.edu.eg'
$ws.Cells.Item(84, 7).Value = 'yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg'
$ws.Cells.Item(85, 7).Value = 'Monica.Eshak@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg'
$ws.Cells.Item(86, 7).Value = 'youstina.magdy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Cells.Item(87, 7).Value = 'eman.samir@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Cells.Item(88, 7).Value = 'maryam.ashraf@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Eman_mohamed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, marinasorial@med.asu.edu.eg'
$ws.Cells.Item(92, 7).Value = 'nahla.nagiub@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, heba@med.asu.edu.eg'
$ws.Cells.Item(93, 7).Value = 'Fatmaelhady@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, AbeerRagheb@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg'
$ws.Cells.Item(96, 7).Value = 'mariam.noureldin@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg'
$ws.Cells.Item(98, 7).Value = 'Madeha.Saeed@med.asu.edu.eg, merna.said@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, maimustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg'
$ws.Cells.Item(104, 7).Value = 'user@user.com, nancy.abdelshafy@med.asu.edu.eg, Amr-Saeed@med.asu.edu.eg'
$ws.Cells.Item(106, 7).Value = 'Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Cells.Item(107, 7).Value = 'yassmen.ahmed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Cells.Item(108, 7).Value = 'yassmen.ahmed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, marinasorial@med.asu.edu.eg'
$ws.Cells.Item(110, 7).Value = 'yassmen.ahmed@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg'
$ws.Cells.Item(111, 7).Value = 'marina_atef@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg'
$ws.Cells.Item(114, 7).Value = 'nahla.nagiub@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, heba@med.asu.edu.eg'
$ws.Cells.Item(115, 7).Value = 'Fatmaelhady@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, AbeerRagheb@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg'
$ws.Cells.Item(118, 7).Value = 'mariam.noureldin@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg'
$ws.Cells.Item(120, 7).Value = 'Madeha.Saeed@med.asu.edu.eg, merna.said@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, maimustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg'
$ws.Cells.Item(126, 7).Value = 'user@user.com, nancy.abdelshafy@med.asu.edu.eg, Amr-Saeed@med.asu.edu.eg'
$ws.Cells.Item(128, 7).Value = 'Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Cells.Item(129, 7).Value = 'yassmen.ahmed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Cells.Item(130, 7).Value = 'yassmen.ahmed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, marinasorial@med.asu.edu.eg'
$ws.Cells.Item(131, 7).Value = 'marina_atef@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Cells.Item(132, 7).Value = 'yassmen.ahmed@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg'
$ws.Cells.Item(133, 7).Value = 'marina_atef@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg'
$ws.Cells.Item(134, 7).Value = 'hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg'
$ws.Cells.Item(137, 7).Value = 'Fatmaelhady@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Nadagouda@med.asu.edu.eg'
$ws.Cells.Item(140, 7).Value = 'aml.awwad@med.asu.edu.eg, aya.saeed@med.asu.edu.eg'
$ws.Cells.Item(142, 7).Value = 'merna.said@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg'
$ws.Cells.Item(144, 7).Value = 'khadija.osama@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg'
$ws.Cells.Item(150, 7).Value = 'yassmen.ahmed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Cells.Item(151, 7).Value = 'yassmen.ahmed@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, marina_atef@med.asu.edu.eg'
$ws.Cells.Item(152, 7).Value = 'marina_atef@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg'
$ws.Cells.Item(153, 7).Value = 'maryam.ashraf@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Eman_mohamed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, marinasorial@med.asu.edu.eg'
$ws.Cells.Item(154, 7).Value = 'Salma.hassan@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg'
$ws.Cells.Item(155, 7).Value = 'marina_atef@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg'
$ws.Cells.Item(156, 7).Value = 'Mohammedeltanany@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm'
$ws.Cells.Item(159, 7).Value = 'Fatmaelhady@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Nadagouda@med.asu.edu.eg'
$ws.Cells.Item(162, 7).Value = 'aml.awwad@med.asu.edu.eg, aya.saeed@med.asu.edu.eg'
$ws.Cells.Item(164, 7).Value = 'merna.said@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg'
$ws.Cells.Item(165, 7).Value = 'Sarah.Mahdy@med.asu.edu.eg, Noran.Mahmoud@med.asu.edu.eg'
$ws.Cells.Item(166, 7).Value = 'khadija.osama@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg'
$ws.Cells.Item(172, 7).Value = 'yassmen.ahmed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Cells.Item(173, 7).Value = 'yassmen.ahmed@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, marina_atef@med.asu.edu.eg'
$ws.Cells.Item(174, 7).Value = 'marina_atef@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg'
$ws.Cells.Item(175, 7).Value = 'maryam.ashraf@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Eman_mohamed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, marinasorial@med.asu.edu.eg'
$ws.Cells.Item(176, 7).Value = 'Salma.hassan@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg'
$ws.Cells.Item(177, 7).Value = 'marina_atef@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg'
